$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - David
$ws.Range("C2").Value = 7
$ws.Range("D2").Value = 5
$ws.Range("G2").Value = 5
$ws.Range("I2").Value = 3
$ws.Range("J2").Value = 5
$ws.Range("M2").Value = 15

# Row 3 - Pedro
$ws.Range("C3").Value = 7
$ws.Range("D3").Value = 3
$ws.Range("G3").Value = 4
$ws.Range("I3").Value = -3
$ws.Range("J3").Value = 3
$ws.Range("M3").Value = 9

# Row 4 - Adonay
$ws.Range("C4").Value = 7
$ws.Range("D4").Value = 7
$ws.Range("G4").Value = 9
$ws.Range("I4").Value = 9
$ws.Range("K4").Value = 2
$ws.Range("M4").Value = 23

# Row 5 - Richard
$ws.Range("C5").Value = 7
$ws.Range("F5").Value = 6
$ws.Range("H5").Value = 7
$ws.Range("I5").Value = -7

# Row 6 - Iván
$ws.Range("C6").Value = 7
$ws.Range("D6").Value = 4
$ws.Range("G6").Value = 8
$ws.Range("I6").Value = 4
$ws.Range("J6").Value = 2
$ws.Range("M6").Value = 16

# Row 7 - Nico
$ws.Range("C7").Value = 7
$ws.Range("F7").Value = 4
$ws.Range("H7").Value = 4
$ws.Range("I7").Value = -2

# Row 8 - Nicolás
$ws.Range("C8").Value = 7
$ws.Range("F8").Value = 5
$ws.Range("H8").Value = 7
$ws.Range("I8").Value = -5

# Row 9 - Vicente
$ws.Range("C9").Value = 7
$ws.Range("F9").Value = 3
$ws.Range("H9").Value = 3
$ws.Range("I9").Value = 1
